$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.856545209884644
$ws.Range("B1").Value = 2.718640565872192
$ws.Range("C1").Value = 2.370614767074585
$ws.Range("D1").Value = 1.585410952568054
$ws.Range("E1").Value = 0.8602614998817444
